$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write a string value to a cell without Excel's "smart" date-entry
# re-interpreting strings that look like dates (e.g. "2025-11-16") as actual
# date serials. Prefixing with a literal leading apostrophe forces text
# entry (the apostrophe itself is NOT stored in the cell), and ClearFormats
# afterwards drops the "quote prefix" cell style that the apostrophe trick
# allocates, so the cell keeps the default (General) style.
# ---------------------------------------------------------------------------
function Set-TextValue($cell, $text) {
    $cell.Value = "'" + $text
}

# ===========================================================================
# Sheet "Chart": append three new daily rows (2025-11-16 .. 2025-11-18)
# ===========================================================================
$chart = $wb.Worksheets.Item("Chart")

Set-TextValue $chart.Cells.Item(44, 1) "2025-11-16"
$chart.Cells.Item(44, 2).Value = 122
$chart.Cells.Item(44, 3).Value = 208
$chart.Cells.Item(44, 4).Value = 32

Set-TextValue $chart.Cells.Item(45, 1) "2025-11-17"
$chart.Cells.Item(45, 2).Value = 122
$chart.Cells.Item(45, 3).Value = 208
$chart.Cells.Item(45, 4).Value = 33

Set-TextValue $chart.Cells.Item(46, 1) "2025-11-18"
$chart.Cells.Item(46, 2).Value = 122
$chart.Cells.Item(46, 3).Value = 208
$chart.Cells.Item(46, 4).Value = 31

$chart.Range("A44:A46").ClearFormats()

# ===========================================================================
# Sheet "Critical issues": refresh the Reason/Source/Validation/Pages table.
# ===========================================================================
$crit = $wb.Worksheets.Item("Critical issues")

$crit.Cells.Item(2, 1).Value = "Not found (404)"
$crit.Cells.Item(2, 2).Value = "Website"
$crit.Cells.Item(2, 3).Value = "Failed"
$crit.Cells.Item(2, 4).Value = 17

$crit.Cells.Item(3, 1).Value = "Duplicate, Google chose different canonical than user"
$crit.Cells.Item(3, 2).Value = "Google systems"
$crit.Cells.Item(3, 3).Value = "Failed"
$crit.Cells.Item(3, 4).Value = 37

$lsq = [char]0x2018
$rsq = [char]0x2019
$crit.Cells.Item(4, 1).Value = "Excluded by " + $lsq + "noindex" + $rsq + " tag"
$crit.Cells.Item(4, 2).Value = "Website"
$crit.Cells.Item(4, 3).Value = "Not Started"
$crit.Cells.Item(4, 4).Value = 15

$crit.Cells.Item(5, 1).Value = "Server error (5xx)"
$crit.Cells.Item(5, 2).Value = "Website"
$crit.Cells.Item(5, 3).Value = "Not Started"
$crit.Cells.Item(5, 4).Value = 2

$crit.Cells.Item(6, 1).Value = "Blocked by robots.txt"
$crit.Cells.Item(6, 2).Value = "Website"
$crit.Cells.Item(6, 3).Value = "Not Started"
$crit.Cells.Item(6, 4).Value = 1

$crit.Cells.Item(7, 1).Value = "Alternate page with proper canonical tag"
$crit.Cells.Item(7, 2).Value = "Website"
$crit.Cells.Item(7, 3).Value = "Started"
$crit.Cells.Item(7, 4).Value = 34

$crit.Cells.Item(8, 1).Value = "Page with redirect"
$crit.Cells.Item(8, 2).Value = "Website"
$crit.Cells.Item(8, 3).Value = "Started"
$crit.Cells.Item(8, 4).Value = 4

$crit.Cells.Item(9, 1).Value = "Crawled - currently not indexed"
$crit.Cells.Item(9, 2).Value = "Google systems"
$crit.Cells.Item(9, 3).Value = "Started"
$crit.Cells.Item(9, 4).Value = 10

$crit.Cells.Item(10, 1).Value = "Discovered - currently not indexed"
$crit.Cells.Item(10, 2).Value = "Google systems"
$crit.Cells.Item(10, 3).Value = "Started"
$crit.Cells.Item(10, 4).Value = 2

Write-Output "done"
